# Apply updated market-price / profit values to the Leve profit tables
# across all 8 crafting-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 2000335.6  # H2: 1667031.6 -> 2000335.6
$ws.Cells.Item(2, 9).Value = 2500284.5  # I2: 2000330 -> 2500284.5
$ws.Cells.Item(2, 11).Value = 2500284.5  # K2: 2000330 -> 2500284.5
$ws.Cells.Item(2, 13).Value = -2500171.5  # M2: -2000217 -> -2500171.5
$ws.Cells.Item(15, 8).Value = 460.51282  # H15: 399.77777 -> 460.51282
$ws.Cells.Item(15, 9).Value = 460.51282  # I15: 399.77777 -> 460.51282
$ws.Cells.Item(15, 11).Value = 1381.53846  # K15: 1199.33331 -> 1381.53846
$ws.Cells.Item(15, 13).Value = -1212.53846  # M15: -1030.33331 -> -1212.53846
$ws.Cells.Item(21, 8).Value = 23663.4  # H21: 25216.75 -> 23663.4
$ws.Cells.Item(21, 9).Value = 21246  # I21: 22511.334 -> 21246
$ws.Cells.Item(21, 11).Value = 21246  # K21: 22511.334 -> 21246
$ws.Cells.Item(21, 13).Value = -20778  # M21: -22043.334 -> -20778
$ws.Cells.Item(23, 8).Value = 23663.4  # H23: 25216.75 -> 23663.4
$ws.Cells.Item(23, 9).Value = 21246  # I23: 22511.334 -> 21246
$ws.Cells.Item(23, 11).Value = 21246  # K23: 22511.334 -> 21246
$ws.Cells.Item(23, 13).Value = -21012  # M23: -22277.334 -> -21012
$ws.Cells.Item(80, 8).Value = 882.5238000000001  # H80: 984 -> 882.5238000000001
$ws.Cells.Item(80, 9).Value = 801.1111  # I80: 870.625 -> 801.1111
$ws.Cells.Item(80, 10).Value = 943.5833  # J80: 1066.4546 -> 943.5833
$ws.Cells.Item(80, 11).Value = 2403.3333  # K80: 2611.875 -> 2403.3333
$ws.Cells.Item(80, 12).Value = 2830.7499  # L80: 3199.3638 -> 2830.7499
$ws.Cells.Item(80, 13).Value = -1405.3333  # M80: -1613.875 -> -1405.3333
$ws.Cells.Item(80, 14).Value = -4826.7499  # N80: -5195.3638 -> -4826.7499
$ws.Cells.Item(83, 8).Value = 882.5238000000001  # H83: 984 -> 882.5238000000001
$ws.Cells.Item(83, 9).Value = 801.1111  # I83: 870.625 -> 801.1111
$ws.Cells.Item(83, 10).Value = 943.5833  # J83: 1066.4546 -> 943.5833
$ws.Cells.Item(83, 11).Value = 7209.9999  # K83: 7835.625 -> 7209.9999
$ws.Cells.Item(83, 12).Value = 8492.2497  # L83: 9598.091400000001 -> 8492.2497
$ws.Cells.Item(83, 13).Value = -2217.9999  # M83: -2843.625 -> -2217.9999
$ws.Cells.Item(83, 14).Value = -18476.2497  # N83: -19582.0914 -> -18476.2497
$ws.Cells.Item(88, 8).Value = 4317.1875  # H88: 3708.2632 -> 4317.1875
$ws.Cells.Item(88, 9).Value = 1900  # I88: 797.3333 -> 1900
$ws.Cells.Item(88, 10).Value = 4478.3335  # J88: 4254.0625 -> 4478.3335
$ws.Cells.Item(88, 11).Value = 1900  # K88: 797.3333 -> 1900
$ws.Cells.Item(88, 12).Value = 4478.3335  # L88: 4254.0625 -> 4478.3335
$ws.Cells.Item(88, 13).Value = -1494  # M88: -391.3333 -> -1494
$ws.Cells.Item(88, 14).Value = -5290.3335  # N88: -5066.0625 -> -5290.3335
$ws.Cells.Item(91, 8).Value = 4317.1875  # H91: 3708.2632 -> 4317.1875
$ws.Cells.Item(91, 9).Value = 1900  # I91: 797.3333 -> 1900
$ws.Cells.Item(91, 10).Value = 4478.3335  # J91: 4254.0625 -> 4478.3335
$ws.Cells.Item(91, 11).Value = 1900  # K91: 797.3333 -> 1900
$ws.Cells.Item(91, 12).Value = 4478.3335  # L91: 4254.0625 -> 4478.3335
$ws.Cells.Item(91, 13).Value = -496  # M91: 606.6667 -> -496
$ws.Cells.Item(91, 14).Value = -7286.3335  # N91: -7062.0625 -> -7286.3335
$ws.Cells.Item(100, 8).Value = 4363.9414  # H100: 3478.2727 -> 4363.9414
$ws.Cells.Item(100, 9).Value = 1688.1111  # I100: 1252 -> 1688.1111
$ws.Cells.Item(100, 11).Value = 1688.1111  # K100: 1252 -> 1688.1111
$ws.Cells.Item(100, 13).Value = -1147.1111  # M100: -711 -> -1147.1111
$ws.Cells.Item(112, 8).Value = 2620.6  # H112: 3035.6667 -> 2620.6
$ws.Cells.Item(112, 10).Value = 3061.1  # J112: 4655.75 -> 3061.1
$ws.Cells.Item(112, 12).Value = 9183.299999999999  # L112: 13967.25 -> 9183.299999999999
$ws.Cells.Item(112, 14).Value = -11399.3  # N112: -16183.25 -> -11399.3
$ws.Cells.Item(114, 8).Value = 10722  # H114: 69999 -> 10722
$ws.Cells.Item(114, 10).Value = 10722  # J114: 69999 -> 10722
$ws.Cells.Item(114, 12).Value = 10722  # L114: 69999 -> 10722
$ws.Cells.Item(114, 14).Value = -19400  # N114: -78677 -> -19400
$ws.Cells.Item(138, 8).Value = 3852.651  # H138: 3988.463 -> 3852.651
$ws.Cells.Item(138, 9).Value = 1953.375  # I138: 2060.2666 -> 1953.375
$ws.Cells.Item(138, 10).Value = 4499.213  # J138: 4730.077 -> 4499.213
$ws.Cells.Item(138, 11).Value = 5860.125  # K138: 6180.7998 -> 5860.125
$ws.Cells.Item(138, 12).Value = 13497.639  # L138: 14190.231 -> 13497.639
$ws.Cells.Item(138, 13).Value = -720.125  # M138: -1040.7998 -> -720.125
$ws.Cells.Item(138, 14).Value = -23777.639  # N138: -24470.231 -> -23777.639

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2444.875  # H2: 3031.8 -> 2444.875
$ws.Cells.Item(2, 9).Value = 2153.2  # I2: 2722 -> 2153.2
$ws.Cells.Item(2, 10).Value = 2931  # J2: 3496.5 -> 2931
$ws.Cells.Item(2, 11).Value = 2153.2  # K2: 2722 -> 2153.2
$ws.Cells.Item(2, 12).Value = 2931  # L2: 3496.5 -> 2931
$ws.Cells.Item(2, 13).Value = -2040.2  # M2: -2609 -> -2040.2
$ws.Cells.Item(2, 14).Value = -3157  # N2: -3722.5 -> -3157
$ws.Cells.Item(61, 8).Value = 3489.8235  # H61: 3568.1516 -> 3489.8235
$ws.Cells.Item(61, 9).Value = 3292.5757  # I61: 3367.1875 -> 3292.5757
$ws.Cells.Item(61, 11).Value = 3292.5757  # K61: 3367.1875 -> 3292.5757
$ws.Cells.Item(61, 13).Value = -3080.5757  # M61: -3155.1875 -> -3080.5757
$ws.Cells.Item(116, 8).Value = 2444.875  # H116: 3031.8 -> 2444.875
$ws.Cells.Item(116, 9).Value = 2153.2  # I116: 2722 -> 2153.2
$ws.Cells.Item(116, 10).Value = 2931  # J116: 3496.5 -> 2931
$ws.Cells.Item(116, 11).Value = 2153.2  # K116: 2722 -> 2153.2
$ws.Cells.Item(116, 12).Value = 2931  # L116: 3496.5 -> 2931
$ws.Cells.Item(116, 13).Value = 140.8000000000002  # M116: -428 -> 140.8000000000002
$ws.Cells.Item(116, 14).Value = -7519  # N116: -8084.5 -> -7519
$ws.Cells.Item(136, 8).Value = 3489.8235  # H136: 3568.1516 -> 3489.8235
$ws.Cells.Item(136, 9).Value = 3292.5757  # I136: 3367.1875 -> 3292.5757
$ws.Cells.Item(136, 11).Value = 9877.7271  # K136: 10101.5625 -> 9877.7271
$ws.Cells.Item(136, 13).Value = -7327.7271  # M136: -7551.5625 -> -7327.7271

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2444.875  # H3: 3031.8 -> 2444.875
$ws.Cells.Item(3, 9).Value = 2153.2  # I3: 2722 -> 2153.2
$ws.Cells.Item(3, 10).Value = 2931  # J3: 3496.5 -> 2931
$ws.Cells.Item(3, 11).Value = 2153.2  # K3: 2722 -> 2153.2
$ws.Cells.Item(3, 12).Value = 2931  # L3: 3496.5 -> 2931
$ws.Cells.Item(3, 13).Value = -2039.2  # M3: -2608 -> -2039.2
$ws.Cells.Item(3, 14).Value = -3159  # N3: -3724.5 -> -3159
$ws.Cells.Item(20, 8).Value = 1844.0555  # H20: 1813.4445 -> 1844.0555
$ws.Cells.Item(20, 9).Value = 1562.625  # I20: 1477.2222 -> 1562.625
$ws.Cells.Item(20, 10).Value = 2069.2  # J20: 2149.6667 -> 2069.2
$ws.Cells.Item(20, 11).Value = 1562.625  # K20: 1477.2222 -> 1562.625
$ws.Cells.Item(20, 12).Value = 2069.2  # L20: 2149.6667 -> 2069.2
$ws.Cells.Item(20, 13).Value = -1315.625  # M20: -1230.2222 -> -1315.625
$ws.Cells.Item(20, 14).Value = -2563.2  # N20: -2643.6667 -> -2563.2
$ws.Cells.Item(105, 8).Value = 4992.5713  # H105: 4606.3125 -> 4992.5713
$ws.Cells.Item(105, 9).Value = 4238.7  # I105: 3849.3333 -> 4238.7
$ws.Cells.Item(105, 11).Value = 4238.7  # K105: 3849.3333 -> 4238.7
$ws.Cells.Item(105, 13).Value = -2491.7  # M105: -2102.3333 -> -2491.7
$ws.Cells.Item(134, 8).Value = 2489.879  # H134: 2449.2354 -> 2489.879
$ws.Cells.Item(134, 9).Value = 2375.7727  # I134: 2320.652 -> 2375.7727
$ws.Cells.Item(134, 11).Value = 7127.3181  # K134: 6961.956 -> 7127.3181
$ws.Cells.Item(134, 13).Value = -4592.3181  # M134: -4426.956 -> -4592.3181

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 904.5  # H107: 947.2 -> 904.5
$ws.Cells.Item(107, 9).Value = 805.4  # I107: 834 -> 805.4
$ws.Cells.Item(107, 11).Value = 805.4  # K107: 834 -> 805.4
$ws.Cells.Item(107, 13).Value = 1114.6  # M107: 1086 -> 1114.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 39999  # H8: 4626.3335 -> 39999
$ws.Cells.Item(8, 9).Value = 39999  # I8: 4626.3335 -> 39999
$ws.Cells.Item(8, 11).Value = 119997  # K8: 13879.0005 -> 119997
$ws.Cells.Item(8, 13).Value = -119858  # M8: -13740.0005 -> -119858
$ws.Cells.Item(12, 8).Value = 89.1875  # H12: 100.21429 -> 89.1875
$ws.Cells.Item(12, 10).Value = 85.5  # J12: 103.875 -> 85.5
$ws.Cells.Item(12, 12).Value = 256.5  # L12: 311.625 -> 256.5
$ws.Cells.Item(12, 14).Value = -602.5  # N12: -657.625 -> -602.5
$ws.Cells.Item(33, 8).Value = 224.6  # H33: 195.8125 -> 224.6
$ws.Cells.Item(33, 9).Value = 234.85715  # I33: 187 -> 234.85715
$ws.Cells.Item(33, 10).Value = 215.625  # J33: 207.14285 -> 215.625
$ws.Cells.Item(33, 11).Value = 1409.1429  # K33: 1122 -> 1409.1429
$ws.Cells.Item(33, 12).Value = 1293.75  # L33: 1242.8571 -> 1293.75
$ws.Cells.Item(33, 13).Value = -1126.1429  # M33: -839 -> -1126.1429
$ws.Cells.Item(33, 14).Value = -1859.75  # N33: -1808.8571 -> -1859.75
$ws.Cells.Item(122, 8).Value = 2912.25  # H122: 1949.8572 -> 2912.25
$ws.Cells.Item(122, 10).Value = 3649.6667  # J122: 2158.1667 -> 3649.6667
$ws.Cells.Item(122, 12).Value = 32847.0003  # L122: 19423.5003 -> 32847.0003
$ws.Cells.Item(122, 14).Value = -37747.0003  # N122: -24323.5003 -> -37747.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10758.444  # H70: 10610.223 -> 10758.444
$ws.Cells.Item(70, 9).Value = 12207.75  # I70: 11299.2 -> 12207.75
$ws.Cells.Item(70, 10).Value = 9599  # J70: 9749 -> 9599
$ws.Cells.Item(70, 11).Value = 12207.75  # K70: 11299.2 -> 12207.75
$ws.Cells.Item(70, 12).Value = 9599  # L70: 9749 -> 9599
$ws.Cells.Item(70, 13).Value = -11937.75  # M70: -11029.2 -> -11937.75
$ws.Cells.Item(70, 14).Value = -10139  # N70: -10289 -> -10139
$ws.Cells.Item(73, 8).Value = 10758.444  # H73: 10610.223 -> 10758.444
$ws.Cells.Item(73, 9).Value = 12207.75  # I73: 11299.2 -> 12207.75
$ws.Cells.Item(73, 10).Value = 9599  # J73: 9749 -> 9599
$ws.Cells.Item(73, 11).Value = 12207.75  # K73: 11299.2 -> 12207.75
$ws.Cells.Item(73, 12).Value = 9599  # L73: 9749 -> 9599
$ws.Cells.Item(73, 13).Value = -11271.75  # M73: -10363.2 -> -11271.75
$ws.Cells.Item(73, 14).Value = -11471  # N73: -11621 -> -11471
$ws.Cells.Item(80, 8).Value = 6285.0713  # H80: 6982.6665 -> 6285.0713
$ws.Cells.Item(80, 9).Value = 3499.3333  # I80: 3949 -> 3499.3333
$ws.Cells.Item(80, 10).Value = 7044.8184  # J80: 7589.4 -> 7044.8184
$ws.Cells.Item(80, 11).Value = 3499.3333  # K80: 3949 -> 3499.3333
$ws.Cells.Item(80, 12).Value = 7044.8184  # L80: 7589.4 -> 7044.8184
$ws.Cells.Item(80, 13).Value = -2501.3333  # M80: -2951 -> -2501.3333
$ws.Cells.Item(80, 14).Value = -9040.8184  # N80: -9585.4 -> -9040.8184
$ws.Cells.Item(83, 8).Value = 6285.0713  # H83: 6982.6665 -> 6285.0713
$ws.Cells.Item(83, 9).Value = 3499.3333  # I83: 3949 -> 3499.3333
$ws.Cells.Item(83, 10).Value = 7044.8184  # J83: 7589.4 -> 7044.8184
$ws.Cells.Item(83, 11).Value = 17496.6665  # K83: 19745 -> 17496.6665
$ws.Cells.Item(83, 12).Value = 35224.092  # L83: 37947 -> 35224.092
$ws.Cells.Item(83, 13).Value = -12504.6665  # M83: -14753 -> -12504.6665
$ws.Cells.Item(83, 14).Value = -45208.092  # N83: -47931 -> -45208.092
$ws.Cells.Item(126, 8).Value = 3219.8  # H126: 3275 -> 3219.8
$ws.Cells.Item(126, 10).Value = 2999.5  # J126: 3000 -> 2999.5
$ws.Cells.Item(126, 12).Value = 8998.5  # L126: 9000 -> 8998.5
$ws.Cells.Item(126, 14).Value = -13938.5  # N126: -13940 -> -13938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 864.8077  # H16: 902.6 -> 864.8077
$ws.Cells.Item(16, 9).Value = 377.1905  # I16: 406.3684 -> 377.1905
$ws.Cells.Item(16, 10).Value = 2912.8  # J16: 2474 -> 2912.8
$ws.Cells.Item(16, 11).Value = 377.1905  # K16: 406.3684 -> 377.1905
$ws.Cells.Item(16, 12).Value = 2912.8  # L16: 2474 -> 2912.8
$ws.Cells.Item(16, 13).Value = -207.1905  # M16: -236.3684 -> -207.1905
$ws.Cells.Item(16, 14).Value = -3252.8  # N16: -2814 -> -3252.8
$ws.Cells.Item(55, 8).Value = 1270.2  # H55: 1377.174 -> 1270.2
$ws.Cells.Item(55, 9).Value = 218.94118  # I55: 230.125 -> 218.94118
$ws.Cells.Item(55, 10).Value = 3504.125  # J55: 3999 -> 3504.125
$ws.Cells.Item(55, 11).Value = 218.94118  # K55: 230.125 -> 218.94118
$ws.Cells.Item(55, 12).Value = 3504.125  # L55: 3999 -> 3504.125
$ws.Cells.Item(55, 13).Value = -45.94118  # M55: -57.125 -> -45.94118
$ws.Cells.Item(55, 14).Value = -3850.125  # N55: -4345 -> -3850.125
$ws.Cells.Item(68, 8).Value = 7538.231  # H68: 4299.9287 -> 7538.231
$ws.Cells.Item(68, 9).Value = 5249.25  # I68: 3187.375 -> 5249.25
$ws.Cells.Item(68, 10).Value = 8555.556  # J68: 4744.95 -> 8555.556
$ws.Cells.Item(68, 11).Value = 5249.25  # K68: 3187.375 -> 5249.25
$ws.Cells.Item(68, 12).Value = 8555.556  # L68: 4744.95 -> 8555.556
$ws.Cells.Item(68, 13).Value = -4500.25  # M68: -2438.375 -> -4500.25
$ws.Cells.Item(68, 14).Value = -10053.556  # N68: -6242.95 -> -10053.556
$ws.Cells.Item(71, 8).Value = 7538.231  # H71: 4299.9287 -> 7538.231
$ws.Cells.Item(71, 9).Value = 5249.25  # I71: 3187.375 -> 5249.25
$ws.Cells.Item(71, 10).Value = 8555.556  # J71: 4744.95 -> 8555.556
$ws.Cells.Item(71, 11).Value = 26246.25  # K71: 15936.875 -> 26246.25
$ws.Cells.Item(71, 12).Value = 42777.78  # L71: 23724.75 -> 42777.78
$ws.Cells.Item(71, 13).Value = -22502.25  # M71: -12192.875 -> -22502.25
$ws.Cells.Item(71, 14).Value = -50265.78  # N71: -31212.75 -> -50265.78
$ws.Cells.Item(132, 8).Value = 1985.849  # H132: 1985.1132 -> 1985.849
$ws.Cells.Item(132, 9).Value = 1867  # I132: 1801.5 -> 1867
$ws.Cells.Item(132, 10).Value = 2024.475  # J132: 2051.0256 -> 2024.475
$ws.Cells.Item(132, 11).Value = 5601  # K132: 5404.5 -> 5601
$ws.Cells.Item(132, 12).Value = 6073.424999999999  # L132: 6153.0768 -> 6073.424999999999
$ws.Cells.Item(132, 13).Value = -3071  # M132: -2874.5 -> -3071
$ws.Cells.Item(132, 14).Value = -11133.425  # N132: -11213.0768 -> -11133.425

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 586.6667  # H14: 574.8570999999999 -> 586.6667
$ws.Cells.Item(14, 10).Value = 0  # J14: 504 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 504 -> 0
$ws.Cells.Item(14, 14).ClearContents()  # N14: -840 -> (blank, since J14/L14 are now 0)
$ws.Cells.Item(126, 8).Value = 3025.611  # H126: 2917.3809 -> 3025.611
$ws.Cells.Item(126, 9).Value = 2566.75  # I126: 2507 -> 2566.75
$ws.Cells.Item(126, 11).Value = 7700.25  # K126: 7521 -> 7700.25
$ws.Cells.Item(126, 13).Value = -5230.25  # M126: -5051 -> -5230.25
$ws.Cells.Item(132, 8).Value = 3046.2  # H132: 1963.7931 -> 3046.2
$ws.Cells.Item(132, 9).Value = 2339.4736  # I132: 1848.2142 -> 2339.4736
$ws.Cells.Item(132, 10).Value = 5284.1665  # J132: 5200 -> 5284.1665
$ws.Cells.Item(132, 11).Value = 7018.4208  # K132: 5544.642599999999 -> 7018.4208
$ws.Cells.Item(132, 12).Value = 15852.4995  # L132: 15600 -> 15852.4995
$ws.Cells.Item(132, 13).Value = -4488.4208  # M132: -3014.642599999999 -> -4488.4208
$ws.Cells.Item(132, 14).Value = -20912.4995  # N132: -20660 -> -20912.4995
